$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01449599266052246
$ws.Range("C2").Value = 0.03331928253173828
$ws.Range("D2").Value = 0.002972698211669922
$ws.Range("E2").Value = 0.01921453475952149
$ws.Range("F2").Value = 0.002491474151611328
$ws.Range("G2").Value = 0.02681756019592285
$ws.Range("H2").Value = 0.005090665817260742
$ws.Range("I2").Value = 0.1278950214385986
$ws.Range("J2").Value = 0.008940410614013673
$ws.Range("K2").Value = 0.02933816909790039
$ws.Range("L2").Value = 0.00439305305480957
$ws.Range("M2").Value = 0.02511696815490723
$ws.Range("B3").Value = 0.005902862548828125
$ws.Range("C3").Value = 0.01815061569213867
$ws.Range("D3").Value = 0.004581928253173828
$ws.Range("E3").Value = 0.01521816253662109
$ws.Range("F3").Value = 0.007916784286499024
$ws.Range("G3").Value = 0.01269001960754395
$ws.Range("H3").Value = 0.012579345703125
$ws.Range("I3").Value = 0.02421746253967285
$ws.Range("J3").Value = 0.008797836303710938
$ws.Range("K3").Value = 0.02429766654968262
$ws.Range("L3").Value = 0.008335208892822266
$ws.Range("M3").Value = 0.01975979804992676
$ws.Range("B4").Value = 0.01591005325317383
$ws.Range("C4").Value = 0.02801060676574707
$ws.Range("D4").Value = 0.01638326644897461
$ws.Range("E4").Value = 0.03671183586120606
$ws.Range("F4").Value = 0.01488223075866699
$ws.Range("G4").Value = 0.03398809432983398
$ws.Range("H4").Value = 0.02066965103149414
$ws.Range("I4").Value = 0.03366847038269043
$ws.Range("J4").Value = 0.01197209358215332
$ws.Range("K4").Value = 0.02116913795471191
$ws.Range("L4").Value = 0.004584217071533203
$ws.Range("M4").Value = 0.01837635040283203
$ws.Range("B5").Value = 0.009705448150634765
$ws.Range("C5").Value = 0.02022037506103515
$ws.Range("D5").Value = 0.01128020286560059
$ws.Range("E5").Value = 0.01858878135681152
$ws.Range("H5").Value = 0.01017975807189941
$ws.Range("I5").Value = 0.02080850601196289
$ws.Range("J5").Value = 0.00831756591796875
$ws.Range("K5").Value = 0.02153096199035644
$ws.Range("B6").Value = 0.03844566345214843
$ws.Range("C6").Value = 0.03099164962768555
$ws.Range("D6").Value = 0.02246980667114258
$ws.Range("E6").Value = 0.02287020683288574
$ws.Range("F6").Value = 0.02323465347290039
$ws.Range("G6").Value = 0.02307519912719726
$ws.Range("H6").Value = 0.03877906799316407
$ws.Range("I6").Value = 0.03765926361083984
$ws.Range("J6").Value = 0.03306832313537598
$ws.Range("K6").Value = 0.03051161766052246
$ws.Range("L6").Value = 0.02558016777038574
$ws.Range("M6").Value = 0.02948288917541504
